# Applies the edits described by the commit:
#   "Plantilla  Azul OK + Juan Lopez + Esteban Mercado"
#
# 1) AB4 (Simon Garcia - Fortalezas): add spaces around the "|" separators.
# 2) AB5 (Juan Lopez - Fortalezas, rich text): "|C" -> "| C" and
#    " stopping attacks..." -> " Stopping attacks..." (capitalize).
# 3) Sheet view: scroll/selection moved from D2/A1 to AB5/Y1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) AB4: plain string, just replace the whole value ---
$ws.Range("AB4").Value = "Solid defensive presence | Leadership and reading of the game | Strong in the air"

# --- 2) AB5: rich text, edit only the affected runs so the other runs'
#     formatting stays untouched ---
$cell = $ws.Range("AB5")

# "|C" (right after "High defensive awareness ") -> "| C"
$idx = $cell.Characters().Text.IndexOf("|Consistently") + 1
$seg = $cell.Characters($idx, 2)
$seg.Text = "| C"

# " stopping attacks and leading the back line. " -> capitalize "Stopping"
$idx2 = $cell.Characters().Text.IndexOf(" stopping attacks") + 2
$seg2 = $cell.Characters($idx2, 1)
$seg2.Text = "S"

# --- 3) Scroll / selection of the sheet view ---
$ws.Application.ActiveWindow.ScrollColumn = 25 # column Y
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("AB5").Select()
